$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column B first (marcel1's name/last name), then row 1 across, then row 2 across
$ws.Range("B1").Value = "marcel1"
$ws.Range("B2").Value = "rodriguez"
$ws.Range("C1").Value = "marcel2"
$ws.Range("D1").Value = "marcel3"
$ws.Range("C2").Value = "papaya"
$ws.Range("D2").Value = "kiwi"

# Age row
$ws.Range("B3").Value = 39
$ws.Range("C3").Value = 32
$ws.Range("D3").Value = 12

# Salary row
$ws.Range("B4").Value = 340009
$ws.Range("B4").NumberFormat = "#,##0"
$ws.Range("C4").Value = 34000
$ws.Range("D4").Value = 10000

$ws.Range("D4").Select()
